# Update LR-pair data: add Inflammatory-Mac cluster and refresh TPM-derived values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 12,20
$data[0,0] = 'FAPs'
$data[0,1] = 'Saa1'
$data[0,2] = 'Fpr2'
$data[0,3] = 'FAPs'
$data[0,4] = [double]"1"
$data[0,5] = [double]"0.3333333333333333"
$data[0,6] = [double]"0.005000333333333333"
$data[0,7] = [double]"0.015001"
$data[0,8] = [double]"0.2140584197833873"
$data[0,9] = [double]"0.2140584197833873"
$data[0,10] = [double]"2"
$data[0,11] = [double]"0.6666666666666666"
$data[0,12] = [double]"1.109174333333333"
$data[0,13] = [double]"3.327523"
$data[0,14] = [double]"0.1199347472980627"
$data[0,15] = [double]"0.1199347472980627"
$data[0,16] = [double]"0.005546241391444444"
$data[0,17] = [double]"0.049916172523"
$data[0,18] = [double]"0.02567304248374319"
$data[0,19] = [double]"0.02567304248374318"

$data[1,0] = 'FAPs'
$data[1,1] = 'Saa1'
$data[1,2] = 'Fpr2'
$data[1,3] = 'Inflammatory-Mac'
$data[1,4] = [double]"1"
$data[1,5] = [double]"0.3333333333333333"
$data[1,6] = [double]"0.005000333333333333"
$data[1,7] = [double]"0.015001"
$data[1,8] = [double]"0.2140584197833873"
$data[1,9] = [double]"0.2140584197833873"
$data[1,10] = [double]"3"
$data[1,11] = [double]"1"
$data[1,12] = [double]"6.33823"
$data[1,13] = [double]"19.01469"
$data[1,14] = [double]"0.6853512477903235"
$data[1,15] = [double]"0.6853512477903234"
$data[1,16] = [double]"0.03169326274333333"
$data[1,17] = [double]"0.2852393646900001"
$data[1,18] = [double]"0.1467052050985694"
$data[1,19] = [double]"0.1467052050985694"

$data[2,0] = 'FAPs'
$data[2,1] = 'Saa1'
$data[2,2] = 'Fpr2'
$data[2,3] = 'MuSCs'
$data[2,4] = [double]"1"
$data[2,5] = [double]"0.3333333333333333"
$data[2,6] = [double]"0.005000333333333333"
$data[2,7] = [double]"0.015001"
$data[2,8] = [double]"0.2140584197833873"
$data[2,9] = [double]"0.2140584197833873"
$data[2,10] = [double]"1"
$data[2,11] = [double]"0.3333333333333333"
$data[2,12] = [double]"0.09159266666666667"
$data[2,13] = [double]"0.274778"
$data[2,14] = [double]"0.009903892472889619"
$data[2,15] = [double]"0.009903892472889617"
$data[2,16] = [double]"0.0004579938642222222"
$data[2,17] = [double]"0.004121944778"
$data[2,18] = [double]"0.002120011572451336"
$data[2,19] = [double]"0.002120011572451336"

$data[3,0] = 'FAPs'
$data[3,1] = 'Saa1'
$data[3,2] = 'Fpr2'
$data[3,3] = 'Resolving-Mac'
$data[3,4] = [double]"1"
$data[3,5] = [double]"0.3333333333333333"
$data[3,6] = [double]"0.005000333333333333"
$data[3,7] = [double]"0.015001"
$data[3,8] = [double]"0.2140584197833873"
$data[3,9] = [double]"0.2140584197833873"
$data[3,10] = [double]"3"
$data[3,11] = [double]"1"
$data[3,12] = [double]"1.709151333333333"
$data[3,13] = [double]"5.127454"
$data[3,14] = [double]"0.1848101124387242"
$data[3,15] = [double]"0.1848101124387242"
$data[3,16] = [double]"0.008546326383777779"
$data[3,17] = [double]"0.07691693745400001"
$data[3,18] = [double]"0.03956016062862345"
$data[3,19] = [double]"0.03956016062862344"

$data[4,0] = 'Inflammatory-Mac'
$data[4,1] = 'Saa1'
$data[4,2] = 'Fpr2'
$data[4,3] = 'FAPs'
$data[4,4] = [double]"1"
$data[4,5] = [double]"0.3333333333333333"
$data[4,6] = [double]"0.000593"
$data[4,7] = [double]"0.001779"
$data[4,8] = [double]"0.02538563621056237"
$data[4,9] = [double]"0.02538563621056237"
$data[4,10] = [double]"2"
$data[4,11] = [double]"0.6666666666666666"
$data[4,12] = [double]"1.109174333333333"
$data[4,13] = [double]"3.327523"
$data[4,14] = [double]"0.1199347472980627"
$data[4,15] = [double]"0.1199347472980627"
$data[4,16] = [double]"0.0006577403796666666"
$data[4,17] = [double]"0.005919663417"
$data[4,18] = [double]"0.003044619863914348"
$data[4,19] = [double]"0.003044619863914347"

$data[5,0] = 'Inflammatory-Mac'
$data[5,1] = 'Saa1'
$data[5,2] = 'Fpr2'
$data[5,3] = 'Inflammatory-Mac'
$data[5,4] = [double]"1"
$data[5,5] = [double]"0.3333333333333333"
$data[5,6] = [double]"0.000593"
$data[5,7] = [double]"0.001779"
$data[5,8] = [double]"0.02538563621056237"
$data[5,9] = [double]"0.02538563621056237"
$data[5,10] = [double]"3"
$data[5,11] = [double]"1"
$data[5,12] = [double]"6.33823"
$data[5,13] = [double]"19.01469"
$data[5,14] = [double]"0.6853512477903235"
$data[5,15] = [double]"0.6853512477903234"
$data[5,16] = [double]"0.00375857039"
$data[5,17] = [double]"0.03382713351"
$data[5,18] = [double]"0.01739807745286014"
$data[5,19] = [double]"0.01739807745286014"

$data[6,0] = 'Inflammatory-Mac'
$data[6,1] = 'Saa1'
$data[6,2] = 'Fpr2'
$data[6,3] = 'MuSCs'
$data[6,4] = [double]"1"
$data[6,5] = [double]"0.3333333333333333"
$data[6,6] = [double]"0.000593"
$data[6,7] = [double]"0.001779"
$data[6,8] = [double]"0.02538563621056237"
$data[6,9] = [double]"0.02538563621056237"
$data[6,10] = [double]"1"
$data[6,11] = [double]"0.3333333333333333"
$data[6,12] = [double]"0.09159266666666667"
$data[6,13] = [double]"0.274778"
$data[6,14] = [double]"0.009903892472889619"
$data[6,15] = [double]"0.009903892472889617"
$data[6,16] = [double]"5.431445133333333E-05"
$data[6,17] = [double]"0.0004888300620000001"
$data[6,18] = [double]"0.0002514166113853028"
$data[6,19] = [double]"0.0002514166113853027"

$data[7,0] = 'Inflammatory-Mac'
$data[7,1] = 'Saa1'
$data[7,2] = 'Fpr2'
$data[7,3] = 'Resolving-Mac'
$data[7,4] = [double]"1"
$data[7,5] = [double]"0.3333333333333333"
$data[7,6] = [double]"0.000593"
$data[7,7] = [double]"0.001779"
$data[7,8] = [double]"0.02538563621056237"
$data[7,9] = [double]"0.02538563621056237"
$data[7,10] = [double]"3"
$data[7,11] = [double]"1"
$data[7,12] = [double]"1.709151333333333"
$data[7,13] = [double]"5.127454"
$data[7,14] = [double]"0.1848101124387242"
$data[7,15] = [double]"0.1848101124387242"
$data[7,16] = [double]"0.001013526740666667"
$data[7,17] = [double]"0.009121740666"
$data[7,18] = [double]"0.004691522282402581"
$data[7,19] = [double]"0.00469152228240258"

$data[8,0] = 'MuSCs'
$data[8,1] = 'Saa1'
$data[8,2] = 'Fpr2'
$data[8,3] = 'FAPs'
$data[8,4] = [double]"3"
$data[8,5] = [double]"1"
$data[8,6] = [double]"0.01776633333333333"
$data[8,7] = [double]"0.053299"
$data[8,8] = [double]"0.7605559440060503"
$data[8,9] = [double]"0.7605559440060503"
$data[8,10] = [double]"2"
$data[8,11] = [double]"0.6666666666666666"
$data[8,12] = [double]"1.109174333333333"
$data[8,13] = [double]"3.327523"
$data[8,14] = [double]"0.1199347472980627"
$data[8,15] = [double]"0.1199347472980627"
$data[8,16] = [double]"0.01970596093077778"
$data[8,17] = [double]"0.177353648377"
$data[8,18] = [double]"0.09121708495040519"
$data[8,19] = [double]"0.09121708495040518"

$data[9,0] = 'MuSCs'
$data[9,1] = 'Saa1'
$data[9,2] = 'Fpr2'
$data[9,3] = 'Inflammatory-Mac'
$data[9,4] = [double]"3"
$data[9,5] = [double]"1"
$data[9,6] = [double]"0.01776633333333333"
$data[9,7] = [double]"0.053299"
$data[9,8] = [double]"0.7605559440060503"
$data[9,9] = [double]"0.7605559440060503"
$data[9,10] = [double]"3"
$data[9,11] = [double]"1"
$data[9,12] = [double]"6.33823"
$data[9,13] = [double]"19.01469"
$data[9,14] = [double]"0.6853512477903235"
$data[9,15] = [double]"0.6853512477903234"
$data[9,16] = [double]"0.1126071069233333"
$data[9,17] = [double]"1.01346396231"
$data[9,18] = [double]"0.5212479652388941"
$data[9,19] = [double]"0.5212479652388939"

$data[10,0] = 'MuSCs'
$data[10,1] = 'Saa1'
$data[10,2] = 'Fpr2'
$data[10,3] = 'MuSCs'
$data[10,4] = [double]"3"
$data[10,5] = [double]"1"
$data[10,6] = [double]"0.01776633333333333"
$data[10,7] = [double]"0.053299"
$data[10,8] = [double]"0.7605559440060503"
$data[10,9] = [double]"0.7605559440060503"
$data[10,10] = [double]"1"
$data[10,11] = [double]"0.3333333333333333"
$data[10,12] = [double]"0.09159266666666667"
$data[10,13] = [double]"0.274778"
$data[10,14] = [double]"0.009903892472889619"
$data[10,15] = [double]"0.009903892472889617"
$data[10,16] = [double]"0.001627265846888889"
$data[10,17] = [double]"0.014645392622"
$data[10,18] = [double]"0.00753246428905298"
$data[10,19] = [double]"0.007532464289052979"

$data[11,0] = 'MuSCs'
$data[11,1] = 'Saa1'
$data[11,2] = 'Fpr2'
$data[11,3] = 'Resolving-Mac'
$data[11,4] = [double]"3"
$data[11,5] = [double]"1"
$data[11,6] = [double]"0.01776633333333333"
$data[11,7] = [double]"0.053299"
$data[11,8] = [double]"0.7605559440060503"
$data[11,9] = [double]"0.7605559440060503"
$data[11,10] = [double]"3"
$data[11,11] = [double]"1"
$data[11,12] = [double]"1.709151333333333"
$data[11,13] = [double]"5.127454"
$data[11,14] = [double]"0.1848101124387242"
$data[11,15] = [double]"0.1848101124387242"
$data[11,16] = [double]"0.03036535230511111"
$data[11,17] = [double]"0.273288170746"
$data[11,18] = [double]"0.1405584295276982"
$data[11,19] = [double]"0.1405584295276982"

$ws.Range("A2:T13").Value = $data
